# Update the "Investment Summary" cost table on slide 8 (Table Placeholder 3)
# - Remove the "Training & Change Management" row
# - Re-purpose the remaining rows/columns into a 3-year cost breakdown table

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$sh = $s.Shapes.Item(3)
$tbl = $sh.Table

# Drop the 5th row ("Training & Change Management") entirely - the table
# shrinks from 7 rows to 6 rows and the graphic frame auto-resizes.
$tbl.Rows.Item(5).Delete()

# Header row
$tbl.Cell(1, 1).Shape.TextFrame.TextRange.Text = "Cost Category"
$tbl.Cell(1, 2).Shape.TextFrame.TextRange.Text = "Year 1"
$tbl.Cell(1, 3).Shape.TextFrame.TextRange.Text = "Year 2"
$tbl.Cell(1, 4).Shape.TextFrame.TextRange.Text = "Year 3"
$tbl.Cell(1, 5).Shape.TextFrame.TextRange.Text = "3-Year Total"

# Row 2: Professional Services
$tbl.Cell(2, 1).Shape.TextFrame.TextRange.Text = "Professional Services"
$tbl.Cell(2, 2).Shape.TextFrame.TextRange.Text = "$364,000"
$tbl.Cell(2, 3).Shape.TextFrame.TextRange.Text = "$0"
$tbl.Cell(2, 4).Shape.TextFrame.TextRange.Text = "$0"
$tbl.Cell(2, 5).Shape.TextFrame.TextRange.Text = "$364,000"

# Row 3: Cloud Infrastructure
$tbl.Cell(3, 1).Shape.TextFrame.TextRange.Text = "Cloud Infrastructure"
$tbl.Cell(3, 2).Shape.TextFrame.TextRange.Text = "$8,914"
$tbl.Cell(3, 3).Shape.TextFrame.TextRange.Text = "$8,914"
$tbl.Cell(3, 4).Shape.TextFrame.TextRange.Text = "$8,914"
$tbl.Cell(3, 5).Shape.TextFrame.TextRange.Text = "$26,741"

# Row 4: Software Licenses & Subscriptions
$tbl.Cell(4, 1).Shape.TextFrame.TextRange.Text = "Software Licenses & Subscriptions"
$tbl.Cell(4, 2).Shape.TextFrame.TextRange.Text = "$7,650"
$tbl.Cell(4, 3).Shape.TextFrame.TextRange.Text = "$7,650"
$tbl.Cell(4, 4).Shape.TextFrame.TextRange.Text = "$7,650"
$tbl.Cell(4, 5).Shape.TextFrame.TextRange.Text = "$22,950"

# Row 5: Support & Maintenance
$tbl.Cell(5, 1).Shape.TextFrame.TextRange.Text = "Support & Maintenance"
$tbl.Cell(5, 2).Shape.TextFrame.TextRange.Text = "$0"
$tbl.Cell(5, 3).Shape.TextFrame.TextRange.Text = "$0"
$tbl.Cell(5, 4).Shape.TextFrame.TextRange.Text = "$0"
$tbl.Cell(5, 5).Shape.TextFrame.TextRange.Text = "$0"

# Row 6: Total
$tbl.Cell(6, 1).Shape.TextFrame.TextRange.Text = "TOTAL SOLUTION INVESTMENT"
$tbl.Cell(6, 2).Shape.TextFrame.TextRange.Text = "$380,564"
$tbl.Cell(6, 3).Shape.TextFrame.TextRange.Text = "$16,564"
$tbl.Cell(6, 4).Shape.TextFrame.TextRange.Text = "$16,564"
$tbl.Cell(6, 5).Shape.TextFrame.TextRange.Text = "$413,691"
